# Finished Week 13 logging
# Update target depth counts for the "H" row on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 131
$wsOff.Range("C2").Value = 82
$wsOff.Range("D2").Value = 40
$wsOff.Range("E2").Value = 21
$wsOff.Range("F2").Value = 6

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 137
$wsDef.Range("C2").Value = 98
$wsDef.Range("D2").Value = 33
$wsDef.Range("E2").Value = 17
